# Auto-generated from OOXML diff: updates cached numeric values
# on the Leve profit-calculation sheets (ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 40: Stuck in the Moment | Horn Glue
$ws.Range("H40").Value = 5666.6665
$ws.Range("J40").Value = 5666.6665
$ws.Range("L40").Value = 5666.6665
$ws.Range("N40").Value = -6016.6665

# ALC row 70: Consecrating Congregation | Holy Water
$ws.Range("H70").Value = 15788.777
$ws.Range("I70").Value = 34699.668
$ws.Range("K70").Value = 104099.004
$ws.Range("M70").Value = -103829.004

# ALC row 73: Curbing the Contagion (L) | Holy Water
$ws.Range("H73").Value = 15788.777
$ws.Range("I73").Value = 34699.668
$ws.Range("K73").Value = 104099.004
$ws.Range("M73").Value = -103163.004

# ALC row 86: Filling in the Blanks | Enchanted Aurum Regis Ink
$ws.Range("H86").Value = 2637.9
$ws.Range("J86").Value = 2698
$ws.Range("L86").Value = 2698
$ws.Range("N86").Value = -4944

# ALC row 89: Ink into Antiquity (L) | Enchanted Aurum Regis Ink
$ws.Range("H89").Value = 2637.9
$ws.Range("J89").Value = 2698
$ws.Range("L89").Value = 13490
$ws.Range("N89").Value = -24722

# ALC row 98: The Dotted Line | Enchanted Durium Ink
$ws.Range("H98").Value = 8900.799999999999
$ws.Range("J98").Value = 33054.5
$ws.Range("L98").Value = 33054.5
$ws.Range("N98").Value = -36050.5

# ALC row 103: Let Loose the Juice | Persimmon Tannin
$ws.Range("H103").Value = 45454860
$ws.Range("J103").Value = 125000376
$ws.Range("L103").Value = 375001128
$ws.Range("N103").Value = -375002300

# ALC row 113: Amaro Kart | Starch Glue
$ws.Range("H113").Value = 5222.222
$ws.Range("J113").Value = 5428.5713
$ws.Range("L113").Value = 5428.5713
$ws.Range("N113").Value = -11936.5713

# ALC row 122: Wishful Inking | Enchanted High Durium Ink
$ws.Range("H122").Value = 8900.799999999999
$ws.Range("J122").Value = 33054.5
$ws.Range("L122").Value = 99163.5
$ws.Range("N122").Value = -104063.5

# ALC row 132: Fast-forwarding Flora | Growth Formula Lambda
$ws.Range("H132").Value = 4265.6855
$ws.Range("I132").Value = 1942.6364
$ws.Range("J132").Value = 8197
$ws.Range("K132").Value = 5827.9092
$ws.Range("L132").Value = 24591
$ws.Range("M132").Value = -3297.9092
$ws.Range("N132").Value = -29651

$ws = $wb.Worksheets.Item("ARM")
# ARM row 122: Haste for High Durium | High Durium Nugget
$ws.Range("H122").Value = 2384.8667
$ws.Range("I122").Value = 2341
$ws.Range("K122").Value = 7023
$ws.Range("M122").Value = -4573

$ws = $wb.Worksheets.Item("BSM")
# BSM row 86: Through Thick and Thin | Adamantite Nugget
$ws.Range("H86").Value = 2928
$ws.Range("I86").Value = 2307.5
$ws.Range("K86").Value = 2307.5
$ws.Range("M86").Value = -1184.5

# BSM row 89: Piercing Eyes Deserve Piercing Shafts (L) | Adamantite Nugget
$ws.Range("H89").Value = 2928
$ws.Range("I89").Value = 2307.5
$ws.Range("K89").Value = 11537.5
$ws.Range("M89").Value = -5921.5

# BSM row 105: Ingot to Wing It | Molybdenum Ingot
$ws.Range("H105").Value = 2635044.2
$ws.Range("I105").Value = 3128614
$ws.Range("K105").Value = 3128614
$ws.Range("M105").Value = -3126867

# BSM row 107: The Gold Experience | Deepgold Nugget
$ws.Range("H107").Value = 5541.9287
$ws.Range("I107").Value = 2972.5417
$ws.Range("J107").Value = 20958.25
$ws.Range("K107").Value = 2972.5417
$ws.Range("L107").Value = 20958.25
$ws.Range("M107").Value = -1052.5417
$ws.Range("N107").Value = -24798.25

$ws = $wb.Worksheets.Item("CRP")
# CRP row 31: Wall Not Found | Walnut Lumber
$ws.Range("H31").Value = 1508.6428
$ws.Range("I31").Value = 1601.9131
$ws.Range("J31").Value = 1079.6
$ws.Range("K31").Value = 1601.9131
$ws.Range("L31").Value = 1079.6
$ws.Range("M31").Value = -1306.9131
$ws.Range("N31").Value = -1669.6

# CRP row 34: Armoires of the Rich and Famous | Walnut Lumber
$ws.Range("H34").Value = 1508.6428
$ws.Range("I34").Value = 1601.9131
$ws.Range("J34").Value = 1079.6
$ws.Range("K34").Value = 1601.9131
$ws.Range("L34").Value = 1079.6
$ws.Range("M34").Value = -1399.9131
$ws.Range("N34").Value = -1483.6

# CRP row 41: The Lone Bowman | Oak Longbow
$ws.Range("H41").Value = 12764.5
$ws.Range("I41").Value = 8775
$ws.Range("J41").Value = 14094.333
$ws.Range("K41").Value = 8775
$ws.Range("L41").Value = 14094.333
$ws.Range("M41").Value = -8347
$ws.Range("N41").Value = -14950.333

# CRP row 50: The Arsenal of Theocracy | Cobalt Halberd
$ws.Range("H50").Value = 11110.556
$ws.Range("J50").Value = 11110.556
$ws.Range("L50").Value = 11110.556
$ws.Range("N50").Value = -12360.556

# CRP row 51: Greenstone for Greenhorns | Jade Crook
$ws.Range("H51").Value = 11249.375
$ws.Range("J51").Value = 11249.375
$ws.Range("L51").Value = 11249.375
$ws.Range("N51").Value = -12721.375

# CRP row 58: You Do the Heavy Lifting | Mahogany Lumber
$ws.Range("H58").Value = 1683.0667
$ws.Range("I58").Value = 1646.2222
$ws.Range("J58").Value = 1738.3334
$ws.Range("K58").Value = 1646.2222
$ws.Range("L58").Value = 1738.3334
$ws.Range("M58").Value = -1443.2222
$ws.Range("N58").Value = -2144.3334

# CRP row 60: Bowing to Greater Power | Yew Longbow
$ws.Range("H60").Value = 10720.429
$ws.Range("J60").Value = 10832.833
$ws.Range("L60").Value = 10832.833
$ws.Range("N60").Value = -11854.833

# CRP row 61: Incant Now, Think Later | Jade Crook
$ws.Range("H61").Value = 11249.375
$ws.Range("J61").Value = 11249.375
$ws.Range("L61").Value = 11249.375
$ws.Range("N61").Value = -11945.375

# CRP row 132: Hull Lotta Damage | Ginseng Lumber
$ws.Range("H132").Value = 3085.9285
$ws.Range("I132").Value = 2938.6924
$ws.Range("K132").Value = 8816.0772
$ws.Range("M132").Value = -6286.0772

# CRP row 136: Turali Quality | Dark Mahogany Lumber
$ws.Range("H136").Value = 1683.0667
$ws.Range("I136").Value = 1646.2222
$ws.Range("J136").Value = 1738.3334
$ws.Range("K136").Value = 4938.6666
$ws.Range("L136").Value = 5215.0002
$ws.Range("M136").Value = -2388.6666
$ws.Range("N136").Value = -10315.0002

$ws = $wb.Worksheets.Item("CUL")
# CUL row 88: Don't Let It Fall Apart | Liver-cheese Sandwich
$ws.Range("H88").Value = 15000
$ws.Range("J88").Value = 15000
$ws.Range("L88").Value = 45000
$ws.Range("N88").Value = -45856

# CUL row 91: Better Come Back with a Sandwich (L) | Liver-cheese Sandwich
$ws.Range("H91").Value = 15000
$ws.Range("J91").Value = 15000
$ws.Range("L91").Value = 45000
$ws.Range("N91").Value = -47964

# CUL row 107: Slippery Service | Frantoio Oil
$ws.Range("H107").Value = 867.8570999999999
$ws.Range("I107").Value = 561.125
$ws.Range("J107").Value = 990.55
$ws.Range("K107").Value = 1683.375
$ws.Range("L107").Value = 2971.65
$ws.Range("M107").Value = 236.625
$ws.Range("N107").Value = -6811.65

# CUL row 137: Creative Chocolate | Gateau au Chocolat
$ws.Range("H137").Value = 5885337
$ws.Range("J137").Value = 7186.6
$ws.Range("L137").Value = 21559.8
$ws.Range("N137").Value = -31759.8

$ws = $wb.Worksheets.Item("GSM")
# GSM row 97: If I'd a Koppranickel for Every Time... | Koppranickel Ingot
$ws.Range("H97").Value = 25281.926
$ws.Range("I97").Value = 33476.65
$ws.Range("J97").Value = 1868.4286
$ws.Range("K97").Value = 33476.65
$ws.Range("L97").Value = 1868.4286
$ws.Range("M97").Value = -32980.65
$ws.Range("N97").Value = -2860.4286

# GSM row 122: Awarding Academic Excellence | Ametrine
$ws.Range("H122").Value = 3215
$ws.Range("I122").Value = 2953.3333
$ws.Range("K122").Value = 8859.999899999999
$ws.Range("M122").Value = -6409.999899999999

# GSM row 132: On Board for Lar | Lar Ingot
$ws.Range("H132").Value = 4251.385
$ws.Range("I132").Value = 4251.385
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 12754.155
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -10224.155
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
# LTW row 40: Best Served Toad | Toad Leather
$ws.Range("H40").Value = 4338.5
$ws.Range("I40").Value = 3923.75
$ws.Range("K40").Value = 3923.75
$ws.Range("M40").Value = -3787.75

# LTW row 122: Hell on Leather | Gaja Leather
$ws.Range("H122").Value = 13949.7
$ws.Range("I122").Value = 20000
$ws.Range("J122").Value = 4874.25
$ws.Range("K122").Value = 60000
$ws.Range("L122").Value = 14622.75
$ws.Range("M122").Value = -57550
$ws.Range("N122").Value = -19522.75

# LTW row 132: Tenets of Tanning | Silver Lobo Leather
$ws.Range("H132").Value = 2837.72
$ws.Range("I132").Value = 2207.45
$ws.Range("K132").Value = 6622.349999999999
$ws.Range("M132").Value = -4092.349999999999

# LTW row 136: Respect for Br'aax | Br'aax Leather
$ws.Range("H136").Value = 2883.1
$ws.Range("I136").Value = 1820.5385
$ws.Range("J136").Value = 4856.4287
$ws.Range("K136").Value = 5461.6155
$ws.Range("L136").Value = 14569.2861
$ws.Range("M136").Value = -2911.6155
$ws.Range("N136").Value = -19669.2861

$ws = $wb.Worksheets.Item("WVR")
# WVR row 81: Where the Dragonflies, the Net Catches | Crawler Silk
$ws.Range("H81").Value = 3198
$ws.Range("I81").Value = 3198
$ws.Range("K81").Value = 6396
$ws.Range("M81").Value = -5335

# WVR row 84: To Kill a Dragon on Nameday (L) | Crawler Silk
$ws.Range("H84").Value = 3198
$ws.Range("I84").Value = 3198
$ws.Range("K84").Value = 31980
$ws.Range("M84").Value = -26676

# WVR row 107: Flax Wax | Bright Linen Yarn
$ws.Range("H107").Value = 71430080
$ws.Range("I107").Value = 1864
$ws.Range("J107").Value = 125001250
$ws.Range("K107").Value = 5592
$ws.Range("L107").Value = 375003750
$ws.Range("M107").Value = -3672
$ws.Range("N107").Value = -375007590

# WVR row 119: A Job Well Done | Dwarven Cotton Gaskins of Fending
$ws.Range("H119").Value = 29898
$ws.Range("J119").Value = 29898
$ws.Range("L119").Value = 29898
$ws.Range("N119").Value = -39574

# WVR row 122: Heavy Armoire | Dark Hempen Cloth
$ws.Range("H122").Value = 2008.4828
$ws.Range("I122").Value = 1526.8636
$ws.Range("J122").Value = 3522.1428
$ws.Range("K122").Value = 4580.5908
$ws.Range("L122").Value = 10566.4284
$ws.Range("M122").Value = -2130.5908
$ws.Range("N122").Value = -15466.4284

# WVR row 132: Comfy Cabins | Snow Cotton Cloth
$ws.Range("H132").Value = 15664.723
$ws.Range("I132").Value = 13690.23
$ws.Range("J132").Value = 20798.4
$ws.Range("K132").Value = 41070.69
$ws.Range("L132").Value = 62395.2
$ws.Range("M132").Value = -38540.69
$ws.Range("N132").Value = -67455.20000000001

# WVR row 136: Weaving the Envelope | Sarcenet Cloth
$ws.Range("H136").Value = 944.8421
$ws.Range("I136").Value = 964
$ws.Range("K136").Value = 2892
$ws.Range("M136").Value = -342
